# Update profit.py after running on 2025-08-25
#
# Sheet1: append a new row (row 8) with the day's date/profit, matching
# the plain-text date style already used by the existing rows (no date
# auto-conversion, no leftover number-format override on the cell).
#
# Sheet2: overwrite the single summary row (row 2) with the refreshed
# date + the two recomputed ratios.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Force text entry for the date so Excel doesn't auto-convert the
# "MM/DD/YYYY" string into a date serial (as it would for a plain
# General-formatted cell), then drop the temporary "@" number format so
# the cell ends up styled exactly like its neighbours (default style).
$ws1.Range("A8").NumberFormat = "@"
$ws1.Range("A8").Value = "08/25/2025"
$ws1.Range("A8").ClearFormats()

$ws1.Range("B8").Value = 14927.67

$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A2").NumberFormat = "@"
$ws2.Range("A2").Value = "08/25/2025"
$ws2.Range("A2").ClearFormats()

$ws2.Range("B2").Value = 0.09453157997486517
$ws2.Range("C2").Value = 0.9054684200251348
